# Add a new slide ("Controls") to the end of the presentation, using the
# same "Title and Content" layout as the other content slides.
$p = $ppt.ActivePresentation
$newIndex = $p.Slides.Count + 1
$s = $p.Slides.Add($newIndex, 2)

# Title placeholder
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Controls"

# Body placeholder - build it run by run so the paragraph/run layout
# mirrors the authored slide.
$body = $s.Shapes.Item(2).TextFrame.TextRange

$r = $body
$body.Text = "For "
$r = $r.InsertAfter("moving")
$r = $r.InsertAfter(": a(")
$r = $r.InsertAfter("left")
$r = $r.InsertAfter("),s(")
$r = $r.InsertAfter("back")
$r = $r.InsertAfter("),d(")
$r = $r.InsertAfter("right")
$r = $r.InsertAfter("),w(")
$r = $r.InsertAfter("forward")
$r = $r.InsertAfter(")")

$r = $r.InsertAfter("`rFor ")
$r = $r.InsertAfter("jumping")
$r = $r.InsertAfter(": ")
$r = $r.InsertAfter("space")
$r = $r.InsertAfter(" ")
$r = $r.InsertAfter("bar")

$r = $r.InsertAfter("`rFor ")
$r = $r.InsertAfter("aiming")
$r = $r.InsertAfter(": ")
$r = $r.InsertAfter("move")
$r = $r.InsertAfter(" ")
$r = $r.InsertAfter("mouse")

$r = $r.InsertAfter("`rFor ")
$r = $r.InsertAfter("shooting")
$r = $r.InsertAfter(": ")
$r = $r.InsertAfter("mouse")
$r = $r.InsertAfter(" ")
$r = $r.InsertAfter("left")
$r = $r.InsertAfter(" ")
$r = $r.InsertAfter("button")
